$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.905.95"
$ws.Range("E2").Value = "  -0.12%  "
$ws.Range("D3").Value = "1.814.37"
$ws.Range("E3").Value = "  +0.18%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "308.97"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.54%  "
$ws.Range("E6").Value = "  +0.06%  "
$ws.Range("E7").Value = "  +0.95%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3657"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -1.49%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07363"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -0.16%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8688"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -0.61%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.23"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -1.21%  "
$ws.Range("D12").Value = "1.831.94"
$ws.Range("E12").Value = "  +2.24%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.383"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +0.42%  "
$ws.Range("E14").Value = "  +0.77%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.501"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -0.37%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "91.07"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -1.41%  "
$ws.Range("E17").Value = "  +0.01%  "
$ws.Range("E18").Value = "  -0.10%  "
$ws.Range("E19").Value = "  +0.07%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.62"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -0.75%  "
$ws.Range("D21").Value = "26.925.44"
$ws.Range("E21").Value = "  -0.09%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.293"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -0.41%  "
$ws.Range("E23").Value = "  -0.80%  "
$ws.Range("D24").Value = "2.047.74"
$ws.Range("E24").Value = "  +1.71%  "
$ws.Range("E25").Value = "  -0.18%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "150.91"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -0.32%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.38"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -0.02%  "
$ws.Range("E28").Value = "  -0.68%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.259"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -0.64%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "115.78"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -0.13%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08910"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -0.14%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.7541"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +0.27%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.163"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +0.51%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.477"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +0.67%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.908"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -0.56%  "
$ws.Range("E36").Value = "  +0.11%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.087"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -1.64%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.05276"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +0.65%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01942"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -1.77%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.971"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +1.47%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "7.183"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -0.12%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.5277"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -0.65%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.320"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -4.44%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.1650"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -0.93%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.423"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -1.09%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4845"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -2.81%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.44"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +0.48%  "
$ws.Range("E48").Value = "  +0.11%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "103.17"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -1.04%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.658"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -0.71%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06287"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -0.18%  "
